$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = (2.525451729410438 / 1000000)
$ws.Range("E2").Value = (2.525451729410438 / 1000000)
$ws.Range("D3").Value = (8.300367348431201 / 10000000)
$ws.Range("E3").Value = (8.300367348431201 / 10000000)
$ws.Range("D4").Value = 0.9999999994448339
$ws.Range("E4").Value = 0.9999999994448339
$ws.Range("D5").Value = 0.9999999947117257
$ws.Range("E5").Value = 0.9999999947117257
$ws.Range("D6").Value = 0.9999999999986409
$ws.Range("E6").Value = 0.9999999999986409
$ws.Range("D8").Value = (59.2151658038622 / 100000000000000000)
$ws.Range("E8").Value = 0.9999999999999994
$ws.Range("D9").Value = (2.646024416517223 / 1000000)
$ws.Range("E9").Value = 0.9999973539755835
$ws.Range("D10").Value = (45.4343339558829 / 100000000)
$ws.Range("E10").Value = 0.9999995456566605
$ws.Range("D11").Value = 0.001655588090924821
$ws.Range("E11").Value = 0.9983444119090752
$ws.Range("F11").Value = 13.66068553924561
